$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order corresponding to data series: C,D,E,F,G,I,J,K,M,N,O (H and L stay 0)
$cols = @(3, 4, 5, 6, 7, 9, 10, 11, 13, 14, 15)

$data = @{}
$data[2] = @(4.136103668755909, 9.012363208363849, 13.327999910634, 32.56003849646859, 3.651790664247824, 22.64048159906616, 9.901962788081462, 13.64782969350226, 16.92893180262911, 18.06803169272159, 24.54088136003697)
$data[3] = @(4.037996911736511, 8.998688632331502, 13.34443554280399, 32.60905162601514, 3.653857433784242, 22.68370602297985, 9.927420582521934, 13.25998805950048, 16.78052238570713, 18.11251829228348, 24.60708639993614)
$data[4] = @(3.975801007628341, 8.991521330507821, 13.35667525222626, 32.64738458371975, 3.655194259313915, 22.71578979570894, 9.944168991613166, 13.01750910867802, 16.69119879043544, 18.1416196970103, 24.65350955624586)
$data[5] = @(3.949983286823394, 8.988911885105663, 13.36220333471253, 32.66507285963042, 3.655756134428396, 22.7302549925714, 9.951275374952798, 12.91775004992152, 16.65528420021464, 18.15392879419481, 24.67387512385417)
$data[6] = @(3.945668352667802, 8.988497452271597, 13.36315390315136, 32.66813469920005, 3.655850468218596, 22.73274082558962, 9.952472382497954, 12.90113236636912, 16.64935087409982, 18.15599991602409, 24.67734413238061)
$data[7] = @(3.97545470613456, 8.991484875361641, 13.35674761812538, 32.64761477006687, 3.655201767618251, 22.71597925172478, 9.944263691481416, 13.01616735856496, 16.69071242527526, 18.14178387851821, 24.65377835689015)
$data[8] = @(4.102694387917117, 9.007394554349263, 13.33322117076609, 32.57522586340565, 3.652489242634989, 22.65423294093094, 9.910508958867922, 13.51508728943929, 16.87740647308716, 18.0830003843609, 24.56250839433621)
$data[9] = @(4.335868514080262, 9.048238602580669, 13.30412200304901, 32.49880662160928, 3.64770566165826, 22.5772641559514, 9.853165587741151, 14.45313062982262, 17.25632390754144, 17.98186437829785, 24.42949962667971)
$data[10] = @(4.496235457655319, 9.083975232151888, 13.29311242227749, 32.48279074678862, 3.64451428982927, 22.54775903794051, 9.81640961203731, 15.11055763065482, 17.54051310101192, 17.91612918699504, 24.36001465299875)
$data[11] = @(4.566632408182397, 9.101440789237026, 13.29034980674922, 32.48423823268269, 3.643131885182669, 22.54023103978595, 9.800851002284642, 15.40143252129763, 17.6706230745141, 17.88807512076518, 24.33457469524115)
$data[12] = @(4.592909048396903, 9.108225024398591, 13.28962588289444, 32.48604215850717, 3.642618323951771, 22.53822887394124, 9.795126121107673, 15.51030546163517, 17.71997525702991, 17.87771694714871, 24.32583103009921)
$data[13] = @(4.587267062218812, 9.10675639517315, 13.28976747477554, 32.48559780981319, 3.642728487898531, 22.53862232813434, 9.79635166095008, 15.48691587918004, 17.70934330725647, 17.87993597588265, 24.32767451590159)
$data[14] = @(4.568801916712391, 9.101995537268373, 13.29028379588331, 32.48436147837668, 3.643089435553101, 22.54004931107346, 9.800376670817581, 15.4104156399567, 17.67468186952773, 17.88721763421584, 24.3338374991239)
$data[15] = @(4.557441456013732, 9.099101462281755, 13.29064199656984, 32.48376771373789, 3.643311817534658, 22.54103390177925, 9.802863824038397, 15.36338826538431, 17.65346036555483, 17.8917123910282, 24.33772847155545)
$data[16] = @(4.491582236162101, 9.082857852758151, 13.29333809715663, 32.48287190407989, 3.644606024821734, 22.54836970833223, 9.817449762181722, 15.09137491140769, 17.53202365038376, 17.91799974335963, 24.36180159822889)
$data[17] = @(4.450515195755115, 9.07320014775795, 13.29556678063288, 32.48455942111256, 3.645417710200197, 22.55438044926754, 9.826695166797801, 14.92233323464492, 17.45771250494118, 17.93459935958274, 24.37815182886074)
$data[18] = @(4.426654795996487, 9.067759283003975, 13.29706008679857, 32.48635210799748, 3.64589110218384, 22.55839241779284, 9.832122254080273, 14.82433904518602, 17.41505112137139, 17.94432111121902, 24.3881365694265)
$data[19] = @(4.418535342099248, 9.06593677913574, 13.29760202680059, 32.48710025756195, 3.646052507968752, 22.55984604186717, 9.833978565508101, 14.79103147043985, 17.40062165634944, 17.94764264890387, 24.39161684545443)
$data[20] = @(4.454911766646747, 9.07421645343023, 13.29530765637032, 32.48429469779454, 3.645330629224685, 22.55368317242089, 9.825699659596234, 14.94040800989899, 17.46561499936264, 17.93281428883954, 24.37635121592013)
$data[21] = @(4.574236025778668, 9.103389319856888, 13.2901234012378, 32.48469054203606, 3.642983147496031, 22.53960713783749, 9.799189902182597, 15.43292094585242, 17.68486083355264, 17.88507163953854, 24.33200310960881)
$data[22] = @(4.649994116230561, 9.123447041461265, 13.28861290656628, 32.49226827665404, 3.641506762705649, 22.53535354947014, 9.782836524560725, 15.74733616509458, 17.82861368085804, 17.8554151118569, 24.30820708852403)
$data[23] = @(4.609768625108247, 9.112652300090785, 13.28924754305456, 32.48755446688423, 3.642289461803993, 22.53717103707638, 9.791475747124419, 15.58024001541409, 17.75185991460458, 17.8711020946296, 24.32043193862972)
$data[24] = @(4.452924855872928, 9.073756634146338, 13.29542414602409, 32.48441181712632, 3.645369977525446, 22.55399667846431, 9.826149380192005, 14.93223891910482, 17.46204208852031, 17.93362076352283, 24.37716345167998)
$data[25] = @(4.274645052279875, 9.036171538558586, 13.31017145435123, 32.51244225201655, 3.648942756201818, 22.59334602224407, 9.867733243571658, 14.20446273037448, 17.15265129861978, 18.0077159576612, 24.46053749629013)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}
